$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in the lookup table: "rootsNPlaintain" -> "rootsNPlantain"
# This cell (G13) is the source of the VLOOKUP used by column D, so
# dependent formula cells will recalc to the corrected spelling automatically.
$ws.Range("G13").Value = "rootsNPlantain"

# Update the view's scroll position / selection to match the saved state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("K9").Select()
